$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2's displayed text (drop trailing slash) before wiring the hyperlink
$ws.Range("B2").Value2 = "https://www.notonthehighstreet.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.notonthehighstreet.com")

# New column E: header + value
$ws.Range("E1").Value2 = "mothers_day_gifts_path"
$ws.Range("E2").Value2 = "/mothers-day/best-mothers-day-gifts"

# Match formatting of the existing header/body cells in column E
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column E width (~29 character units)
$ws.Columns.Item(5).ColumnWidth = 28.17

# Selection moves to E4
[void]$ws.Range("E4").Select()
